# Refresh market-board derived profit figures (currentAveragePrice* / LevePrice* / LeveProfit*)
# across all Leve sheets, per the scheduled price-data sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 606.75
$ws.Range("I2").Value = 380.16666
$ws.Range("J2").Value = 833.3333
$ws.Range("K2").Value = 380.16666
$ws.Range("L2").Value = 833.3333
$ws.Range("M2").Value = -267.16666
$ws.Range("N2").Value = -1059.3333
# Row 106 (Leve Item ID 19903)
$ws.Range("H106").Value = 3717.7778
$ws.Range("I106").Value = 4005.2632
$ws.Range("J106").Value = 3507.6924
$ws.Range("K106").Value = 4005.2632
$ws.Range("L106").Value = 3507.6924
$ws.Range("M106").Value = -3374.2632
$ws.Range("N106").Value = -4769.6924
# Row 111 (Leve Item ID 27768)
$ws.Range("H111").Value = 3535.5715
$ws.Range("I111").Value = 2458.1667
$ws.Range("K111").Value = 7374.500100000001
$ws.Range("M111").Value = -4307.500100000001
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 29241602
$ws.Range("I132").Value = 5556507
$ws.Range("J132").Value = 55558372
$ws.Range("K132").Value = 16669521
$ws.Range("L132").Value = 166675116
$ws.Range("M132").Value = -16666991
$ws.Range("N132").Value = -166680176
# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 897751.5600000001
$ws.Range("I137").Value = 2036.4
$ws.Range("J137").Value = 1324282.6
$ws.Range("K137").Value = 6109.200000000001
$ws.Range("L137").Value = 3972847.8
$ws.Range("M137").Value = -3559.200000000001
$ws.Range("N137").Value = -3977947.8

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 3771.2576
$ws.Range("I32").Value = 2595.2097
$ws.Range("J32").Value = 22000
$ws.Range("K32").Value = 2595.2097
$ws.Range("L32").Value = 22000
$ws.Range("M32").Value = -2308.2097
$ws.Range("N32").Value = -22574
# Row 45 (Leve Item ID 27714)
$ws.Range("H45").Value = 7261.1113
$ws.Range("I45").Value = 5590
$ws.Range("J45").Value = 9350
$ws.Range("K45").Value = 5590
$ws.Range("L45").Value = 9350
$ws.Range("M45").Value = -5213
$ws.Range("N45").Value = -10104
# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 1375.4166
$ws.Range("I110").Value = 1428.1428
$ws.Range("J110").Value = 1006.3333
$ws.Range("K110").Value = 1428.1428
$ws.Range("L110").Value = 1006.3333
$ws.Range("M110").Value = 616.8571999999999
$ws.Range("N110").Value = -5096.3333
# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1745.5349
$ws.Range("I132").Value = 1409.6111
$ws.Range("J132").Value = 3473.1428
$ws.Range("K132").Value = 4228.8333
$ws.Range("L132").Value = 10419.4284
$ws.Range("M132").Value = -1698.8333
$ws.Range("N132").Value = -15479.4284

$ws = $wb.Worksheets.Item("BSM")
# Row 105 (Leve Item ID 19947)
$ws.Range("H105").Value = 2515.4849
$ws.Range("I105").Value = 2500
$ws.Range("J105").Value = 2542.5833
$ws.Range("K105").Value = 2500
$ws.Range("L105").Value = 2542.5833
$ws.Range("M105").Value = -753
$ws.Range("N105").Value = -6036.5833

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 10001302
$ws.Range("I31").Value = 776.4
$ws.Range("J31").Value = 20001828
$ws.Range("K31").Value = 776.4
$ws.Range("L31").Value = 20001828
$ws.Range("M31").Value = -481.4
$ws.Range("N31").Value = -20002418
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 10001302
$ws.Range("I34").Value = 776.4
$ws.Range("J34").Value = 20001828
$ws.Range("K34").Value = 776.4
$ws.Range("L34").Value = 20001828
$ws.Range("M34").Value = -574.4
$ws.Range("N34").Value = -20002232
# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 1937
$ws.Range("I105").Value = 1800
$ws.Range("J105").Value = 2005.5
$ws.Range("K105").Value = 1800
$ws.Range("L105").Value = 2005.5
$ws.Range("M105").Value = -53
$ws.Range("N105").Value = -5499.5
# Row 107 (Leve Item ID 27689)
$ws.Range("H107").Value = 1395.6522
$ws.Range("I107").Value = 638.5833
$ws.Range("J107").Value = 2221.5454
$ws.Range("K107").Value = 638.5833
$ws.Range("L107").Value = 2221.5454
$ws.Range("M107").Value = 1281.4167
$ws.Range("N107").Value = -6061.5454
# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 5998.625
$ws.Range("I132").Value = 6455.5
$ws.Range("J132").Value = 4628
$ws.Range("K132").Value = 19366.5
$ws.Range("L132").Value = 13884
$ws.Range("M132").Value = -16836.5
$ws.Range("N132").Value = -18944
# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 18645226
$ws.Range("I134").Value = 2175004.8
$ws.Range("J134").Value = 76924470
$ws.Range("K134").Value = 6525014.399999999
$ws.Range("L134").Value = 230773410
$ws.Range("M134").Value = -6522479.399999999
$ws.Range("N134").Value = -230778480

$ws = $wb.Worksheets.Item("CUL")
# Row 32 (Leve Item ID 4731)
$ws.Range("H32").Value = 14766.556
$ws.Range("I32").Value = 1000
$ws.Range("J32").Value = 18699.857
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 56099.571
$ws.Range("M32").Value = -2717
$ws.Range("N32").Value = -56665.571
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 1754.3214
$ws.Range("J68").Value = 2183.4722
$ws.Range("L68").Value = 6550.4166
$ws.Range("N68").Value = -8172.4166
# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 1754.3214
$ws.Range("J71").Value = 2183.4722
$ws.Range("L71").Value = 19651.2498
$ws.Range("N71").Value = -27763.2498
# Row 107 (Leve Item ID 27838)
$ws.Range("H107").Value = 578616.0600000001
$ws.Range("I107").Value = 665.4761999999999
$ws.Range("J107").Value = 803374.5600000001
$ws.Range("K107").Value = 1996.4286
$ws.Range("L107").Value = 2410123.68
$ws.Range("M107").Value = -76.42859999999973
$ws.Range("N107").Value = -2413963.68
# Row 118 (Leve Item ID 27872)
$ws.Range("H118").Value = 4552.6313
$ws.Range("I118").Value = 950
$ws.Range("K118").Value = 2850
$ws.Range("M118").Value = -1607

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 130.35294
$ws.Range("I2").Value = 14.333333
$ws.Range("J2").Value = 193.63637
$ws.Range("K2").Value = 14.333333
$ws.Range("L2").Value = 193.63637
$ws.Range("M2").Value = 98.666667
$ws.Range("N2").Value = -419.63637
# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1789.1428
$ws.Range("I113").Value = 1543.5
$ws.Range("J113").Value = 2116.6667
$ws.Range("K113").Value = 1543.5
$ws.Range("L113").Value = 2116.6667
$ws.Range("M113").Value = 626.5
$ws.Range("N113").Value = -6456.6667

$ws = $wb.Worksheets.Item("LTW")
# Row 41 (Leve Item ID 3611)
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 38 (Leve Item ID 27990)
$ws.Range("H38").Value = 4000
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 4000
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 4000
$ws.Range("M38").ClearContents()
$ws.Range("N38").Value = -4946
# Row 49 (Leve Item ID 3397)
$ws.Range("H49").Value = 6015.25
$ws.Range("I49").Value = 3999
$ws.Range("J49").Value = 6687.3335
$ws.Range("K49").Value = 3999
$ws.Range("L49").Value = 6687.3335
$ws.Range("N49").Value = -7147.3335
$ws.Range("M49").Value = -3769
# Row 107 (Leve Item ID 27746)
$ws.Range("H107").Value = 15001.429
$ws.Range("I107").Value = 20700.8
$ws.Range("J107").Value = 753
$ws.Range("K107").Value = 62102.39999999999
$ws.Range("L107").Value = 2259
$ws.Range("M107").Value = -60182.39999999999
$ws.Range("N107").Value = -6099
# Row 113 (Leve Item ID 27752)
$ws.Range("H113").Value = 575.93335
$ws.Range("I113").Value = 413.9
$ws.Range("K113").Value = 1241.7
$ws.Range("M113").Value = 928.3000000000002
# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 2765.1714
$ws.Range("I136").Value = 2870.261
$ws.Range("J136").Value = 2563.75
$ws.Range("K136").Value = 8610.782999999999
$ws.Range("L136").Value = 7691.25
$ws.Range("M136").Value = -6060.782999999999
$ws.Range("N136").Value = -12791.25
